# Commit: "update scripts wuth new tpm"
# Refresh NATMI ligand-receptor edge-weight statistics (Pdpn -> Clec1b)
# for the YoungD2 / lrc2p TPM run: columns E-J (ligand-expressing stats),
# M-P (receptor-expressing stats) and Q-T (edge weights / specificity)
# are recomputed downstream of a new TPM table and overwritten per row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.698400333333333
$ws.Range("H2").Value = 11.095201
$ws.Range("I2").Value = 0.01460277526461191
$ws.Range("J2").Value = 0.01460277526461191
$ws.Range("M2").Value = 2.801728333333333
$ws.Range("N2").Value = 8.405185
$ws.Range("O2").Value = 0.05711011770608418
$ws.Range("P2").Value = 0.05711011770608417
$ws.Range("Q2").Value = 10.36191300190944
$ws.Range("R2").Value = 93.25721701718498
$ws.Range("S2").Value = 0.0008339662141974805
$ws.Range("T2").Value = 0.0008339662141974803
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.698400333333333
$ws.Range("H3").Value = 11.095201
$ws.Range("I3").Value = 0.01460277526461191
$ws.Range("J3").Value = 0.01460277526461191
$ws.Range("M3").Value = 1.469855
$ws.Range("N3").Value = 4.409565
$ws.Range("O3").Value = 0.02996136030112711
$ws.Range("P3").Value = 0.02996136030112711
$ws.Range("Q3").Value = 5.436112221951666
$ws.Range("R3").Value = 48.92500999756499
$ws.Range("S3").Value = 0.0004375190110994241
$ws.Range("T3").Value = 0.000437519011099424
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.698400333333333
$ws.Range("H4").Value = 11.095201
$ws.Range("I4").Value = 0.01460277526461191
$ws.Range("J4").Value = 0.01460277526461191
$ws.Range("M4").Value = 44.27998366666666
$ws.Range("N4").Value = 132.839951
$ws.Range("O4").Value = 0.9025982459256344
$ws.Range("P4").Value = 0.9025982459256343
$ws.Range("Q4").Value = 163.7651063527945
$ws.Range("R4").Value = 1473.885957175151
$ws.Range("S4").Value = 0.01318043933948495
$ws.Range("T4").Value = 0.01318043933948495
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.698400333333333
$ws.Range("H5").Value = 11.095201
$ws.Range("I5").Value = 0.01460277526461191
$ws.Range("J5").Value = 0.01460277526461191
$ws.Range("M5").Value = 0.5067863333333333
$ws.Range("N5").Value = 1.520359
$ws.Range("O5").Value = 0.01033027606715431
$ws.Range("P5").Value = 0.01033027606715431
$ws.Range("Q5").Value = 1.874298744128778
$ws.Range("R5").Value = 16.868688697159
$ws.Range("S5").Value = 0.0001508506998300534
$ws.Range("T5").Value = 0.0001508506998300533
# Row 6
$ws.Range("I6").Value = 0.4072092715200794
$ws.Range("J6").Value = 0.4072092715200794
$ws.Range("M6").Value = 2.801728333333333
$ws.Range("N6").Value = 8.405185
$ws.Range("O6").Value = 0.05711011770608418
$ws.Range("P6").Value = 0.05711011770608417
$ws.Range("Q6").Value = 288.9496666628405
$ws.Range("R6").Value = 2600.546999965565
$ws.Range("S6").Value = 0.02325576942752052
$ws.Range("T6").Value = 0.02325576942752052
# Row 7
$ws.Range("I7").Value = 0.4072092715200794
$ws.Range("J7").Value = 0.4072092715200794
$ws.Range("M7").Value = 1.469855
$ws.Range("N7").Value = 4.409565
$ws.Range("O7").Value = 0.02996136030112711
$ws.Range("P7").Value = 0.02996136030112711
$ws.Range("Q7").Value = 151.5900407757983
$ws.Range("R7").Value = 1364.310366982185
$ws.Range("S7").Value = 0.0122005437019726
$ws.Range("T7").Value = 0.01220054370197259
# Row 8
$ws.Range("I8").Value = 0.4072092715200794
$ws.Range("J8").Value = 0.4072092715200794
$ws.Range("M8").Value = 44.27998366666666
$ws.Range("N8").Value = 132.839951
$ws.Range("O8").Value = 0.9025982459256344
$ws.Range("P8").Value = 0.9025982459256343
$ws.Range("Q8").Value = 4566.712042740055
$ws.Range("R8").Value = 41100.40838466049
$ws.Range("S8").Value = 0.367546374198679
$ws.Range("T8").Value = 0.367546374198679
# Row 9
$ws.Range("I9").Value = 0.4072092715200794
$ws.Range("J9").Value = 0.4072092715200794
$ws.Range("M9").Value = 0.5067863333333333
$ws.Range("N9").Value = 1.520359
$ws.Range("O9").Value = 0.01033027606715431
$ws.Range("P9").Value = 0.01033027606715431
$ws.Range("Q9").Value = 52.26621737152122
$ws.Range("R9").Value = 470.395956343691
$ws.Range("S9").Value = 0.004206584191907219
$ws.Range("T9").Value = 0.004206584191907218
# Row 10
$ws.Range("G10").Value = 19.93824633333333
$ws.Range("H10").Value = 59.814739
$ws.Range("I10").Value = 0.07872423321834521
$ws.Range("J10").Value = 0.07872423321834522
$ws.Range("M10").Value = 2.801728333333333
$ws.Range("N10").Value = 8.405185
$ws.Range("O10").Value = 0.05711011770608418
$ws.Range("P10").Value = 0.05711011770608417
$ws.Range("Q10").Value = 55.86154966907944
$ws.Range("R10").Value = 502.7539470217149
$ws.Range("S10").Value = 0.004495950225420917
$ws.Range("T10").Value = 0.004495950225420917
# Row 11
$ws.Range("G11").Value = 19.93824633333333
$ws.Range("H11").Value = 59.814739
$ws.Range("I11").Value = 0.07872423321834521
$ws.Range("J11").Value = 0.07872423321834522
$ws.Range("M11").Value = 1.469855
$ws.Range("N11").Value = 4.409565
$ws.Range("O11").Value = 0.02996136030112711
$ws.Range("P11").Value = 0.02996136030112711
$ws.Range("Q11").Value = 29.30633106428166
$ws.Range("R11").Value = 263.756979578535
$ws.Range("S11").Value = 0.002358685115884801
$ws.Range("T11").Value = 0.0023586851158848
# Row 12
$ws.Range("G12").Value = 19.93824633333333
$ws.Range("H12").Value = 59.814739
$ws.Range("I12").Value = 0.07872423321834521
$ws.Range("J12").Value = 0.07872423321834522
$ws.Range("M12").Value = 44.27998366666666
$ws.Range("N12").Value = 132.839951
$ws.Range("O12").Value = 0.9025982459256344
$ws.Range("P12").Value = 0.9025982459256343
$ws.Range("Q12").Value = 882.8652219819763
$ws.Range("R12").Value = 7945.786997837788
$ws.Range("S12").Value = 0.07105635481471895
$ws.Range("T12").Value = 0.07105635481471895
# Row 13
$ws.Range("G13").Value = 19.93824633333333
$ws.Range("H13").Value = 59.814739
$ws.Range("I13").Value = 0.07872423321834521
$ws.Range("J13").Value = 0.07872423321834522
$ws.Range("M13").Value = 0.5067863333333333
$ws.Range("N13").Value = 1.520359
$ws.Range("O13").Value = 0.01033027606715431
$ws.Range("P13").Value = 0.01033027606715431
$ws.Range("Q13").Value = 10.10443075236678
$ws.Range("R13").Value = 90.939876771301
$ws.Range("S13").Value = 0.0008132430623205461
$ws.Range("T13").Value = 0.0008132430623205461
# Row 14
$ws.Range("G14").Value = 32.263448
$ws.Range("H14").Value = 96.790344
$ws.Range("I14").Value = 0.1273890974319868
$ws.Range("J14").Value = 0.1273890974319868
$ws.Range("M14").Value = 2.801728333333333
$ws.Range("N14").Value = 8.405185
$ws.Range("O14").Value = 0.05711011770608418
$ws.Range("P14").Value = 0.05711011770608417
$ws.Range("Q14").Value = 90.39341639262668
$ws.Range("R14").Value = 813.54074753364
$ws.Range("S14").Value = 0.007275206348812593
$ws.Range("T14").Value = 0.007275206348812591
# Row 15
$ws.Range("G15").Value = 32.263448
$ws.Range("H15").Value = 96.790344
$ws.Range("I15").Value = 0.1273890974319868
$ws.Range("J15").Value = 0.1273890974319868
$ws.Range("M15").Value = 1.469855
$ws.Range("N15").Value = 4.409565
$ws.Range("O15").Value = 0.02996136030112711
$ws.Range("P15").Value = 0.02996136030112711
$ws.Range("Q15").Value = 47.42259036004
$ws.Range("R15").Value = 426.80331324036
$ws.Range("S15").Value = 0.003816750646595144
$ws.Range("T15").Value = 0.003816750646595143
# Row 16
$ws.Range("G16").Value = 32.263448
$ws.Range("H16").Value = 96.790344
$ws.Range("I16").Value = 0.1273890974319868
$ws.Range("J16").Value = 0.1273890974319868
$ws.Range("M16").Value = 44.27998366666666
$ws.Range("N16").Value = 132.839951
$ws.Range("O16").Value = 0.9025982459256344
$ws.Range("P16").Value = 0.9025982459256343
$ws.Range("Q16").Value = 1428.624950470349
$ws.Range("R16").Value = 12857.62455423314
$ws.Range("S16").Value = 0.114981175892161
$ws.Range("T16").Value = 0.114981175892161
# Row 17
$ws.Range("G17").Value = 32.263448
$ws.Range("H17").Value = 96.790344
$ws.Range("I17").Value = 0.1273890974319868
$ws.Range("J17").Value = 0.1273890974319868
$ws.Range("M17").Value = 0.5067863333333333
$ws.Range("N17").Value = 1.520359
$ws.Range("O17").Value = 0.01033027606715431
$ws.Range("P17").Value = 0.01033027606715431
$ws.Range("Q17").Value = 16.35067451261067
$ws.Range("R17").Value = 147.156070613496
$ws.Range("S17").Value = 0.001315964544418043
$ws.Range("T17").Value = 0.001315964544418042
# Row 18
$ws.Range("G18").Value = 15.70300266666667
$ws.Range("H18").Value = 47.109008
$ws.Range("I18").Value = 0.06200178408329911
$ws.Range("J18").Value = 0.06200178408329911
$ws.Range("M18").Value = 2.801728333333333
$ws.Range("N18").Value = 8.405185
$ws.Range("O18").Value = 0.05711011770608418
$ws.Range("P18").Value = 0.05711011770608417
$ws.Range("Q18").Value = 43.99554748960889
$ws.Range("R18").Value = 395.95992740648
$ws.Range("S18").Value = 0.003540929186984429
$ws.Range("T18").Value = 0.003540929186984428
# Row 19
$ws.Range("G19").Value = 15.70300266666667
$ws.Range("H19").Value = 47.109008
$ws.Range("I19").Value = 0.06200178408329911
$ws.Range("J19").Value = 0.06200178408329911
$ws.Range("M19").Value = 1.469855
$ws.Range("N19").Value = 4.409565
$ws.Range("O19").Value = 0.02996136030112711
$ws.Range("P19").Value = 0.02996136030112711
$ws.Range("Q19").Value = 23.08113698461333
$ws.Range("R19").Value = 207.73023286152
$ws.Range("S19").Value = 0.001857657792232413
$ws.Range("T19").Value = 0.001857657792232412
# Row 20
$ws.Range("G20").Value = 15.70300266666667
$ws.Range("H20").Value = 47.109008
$ws.Range("I20").Value = 0.06200178408329911
$ws.Range("J20").Value = 0.06200178408329911
$ws.Range("M20").Value = 44.27998366666666
$ws.Range("N20").Value = 132.839951
$ws.Range("O20").Value = 0.9025982459256344
$ws.Range("P20").Value = 0.9025982459256343
$ws.Range("Q20").Value = 695.3287015976231
$ws.Range("R20").Value = 6257.958314378608
$ws.Range("S20").Value = 0.0559627015578457
$ws.Range("T20").Value = 0.0559627015578457
# Row 21
$ws.Range("G21").Value = 15.70300266666667
$ws.Range("H21").Value = 47.109008
$ws.Range("I21").Value = 0.06200178408329911
$ws.Range("J21").Value = 0.06200178408329911
$ws.Range("M21").Value = 0.5067863333333333
$ws.Range("N21").Value = 1.520359
$ws.Range("O21").Value = 0.01033027606715431
$ws.Range("P21").Value = 0.01033027606715431
$ws.Range("Q21").Value = 7.958067143763556
$ws.Range("R21").Value = 71.622604293872
$ws.Range("S21").Value = 0.0006404955462365741
$ws.Range("T21").Value = 0.000640495546236574
# Row 22
$ws.Range("G22").Value = 78.53120166666668
$ws.Range("H22").Value = 235.593605
$ws.Range("I22").Value = 0.3100728384816776
$ws.Range("J22").Value = 0.3100728384816777
$ws.Range("M22").Value = 2.801728333333333
$ws.Range("N22").Value = 8.405185
$ws.Range("O22").Value = 0.05711011770608418
$ws.Range("P22").Value = 0.05711011770608417
$ws.Range("Q22").Value = 220.0230927602139
$ws.Range("R22").Value = 1980.207834841925
$ws.Range("S22").Value = 0.01770829630314824
$ws.Range("T22").Value = 0.01770829630314823
# Row 23
$ws.Range("G23").Value = 78.53120166666668
$ws.Range("H23").Value = 235.593605
$ws.Range("I23").Value = 0.3100728384816776
$ws.Range("J23").Value = 0.3100728384816777
$ws.Range("M23").Value = 1.469855
$ws.Range("N23").Value = 4.409565
$ws.Range("O23").Value = 0.02996136030112711
$ws.Range("P23").Value = 0.02996136030112711
$ws.Range("Q23").Value = 115.4294794257583
$ws.Range("R23").Value = 1038.865314831825
$ws.Range("S23").Value = 0.009290204033342735
$ws.Range("T23").Value = 0.009290204033342735
# Row 24
$ws.Range("G24").Value = 78.53120166666668
$ws.Range("H24").Value = 235.593605
$ws.Range("I24").Value = 0.3100728384816776
$ws.Range("J24").Value = 0.3100728384816777
$ws.Range("M24").Value = 44.27998366666666
$ws.Range("N24").Value = 132.839951
$ws.Range("O24").Value = 0.9025982459256344
$ws.Range("P24").Value = 0.9025982459256343
$ws.Range("Q24").Value = 3477.360327123706
$ws.Range("R24").Value = 31296.24294411335
$ws.Range("S24").Value = 0.2798712001227447
$ws.Range("T24").Value = 0.2798712001227448
# Row 25
$ws.Range("G25").Value = 78.53120166666668
$ws.Range("H25").Value = 235.593605
$ws.Range("I25").Value = 0.3100728384816776
$ws.Range("J25").Value = 0.3100728384816777
$ws.Range("M25").Value = 0.5067863333333333
$ws.Range("N25").Value = 1.520359
$ws.Range("O25").Value = 0.01033027606715431
$ws.Range("P25").Value = 0.01033027606715431
$ws.Range("Q25").Value = 39.79853974491056
$ws.Range("R25").Value = 358.1868577041951
$ws.Range("S25").Value = 0.003203138022441879
$ws.Range("T25").Value = 0.003203138022441879
